# The template's field code "{ m:'doc.html'.fromHTMLURI() }" (a real Word
# field built from w:fldChar/w:instrText runs) must become literal, plain
# text runs spelling out the same characters (so the parser's new
# TokenIteratorFieldRewriterSplit can tokenize the braces itself instead of
# relying on Word's field machinery). The w:bookmarkStart/w:bookmarkEnd for
# "_GoBack" sitting between "doc.html" and "'.fromHTMLURI()" is preserved.

$d = $word.ActiveDocument

$f = $d.Fields.Item(1)

# Remember where the field starts so we can drop the replacement text in
# its place; deleting the field removes its bookmark too, so we rebuild it
# from the literal OOXML below.
$insertAt = $f.Code.Start - 1
$f.Delete()

$target = $d.Range($insertAt, $insertAt)

$runsXml = '<w:r><w:t>{</w:t></w:r>' +
           '<w:r><w:t>m</w:t></w:r>' +
           '<w:r><w:t>:</w:t></w:r>' +
           '<w:r><w:t>' + "'" + '</w:t></w:r>' +
           '<w:r><w:t>doc.html</w:t></w:r>' +
           '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
           '<w:bookmarkEnd w:id="0"/>' +
           '<w:r><w:t>' + "'" + '.fromHTMLURI()</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
              '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
              '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +
              '</w:document>' +
              '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml)
